# Insert a new data row before row 11 (shifts existing rows 11..65 down to 12..66)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = 7
$ws.Range("B11").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C11").Value = "Ñuble"
$ws.Range("D11").Value = 44613
$ws.Range("E11").Value = 16
$ws.Range("F11").Value = 100112022
$ws.Range("G11").Value = "Arveja Verde"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 60
$ws.Range("K11").Value = 24000
$ws.Range("L11").Value = 25000
$ws.Range("M11").Value = 24500
$ws.Range("N11").Value = "$/saco 25 kilos"
$ws.Range("O11").Value = "Provincia de Diguillín"
$ws.Range("P11").Value = 980
$ws.Range("Q11").Value = 25
$ws.Range("R11").Value = "Hortaliza"
